$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 1959.9565
$ws.Cells.Item(40, 10).Value = 2089.6206
$ws.Cells.Item(40, 12).Value = 2089.6206
$ws.Cells.Item(40, 14).Value = -2439.6206
$ws.Cells.Item(112, 8).Value = 1512.421
$ws.Cells.Item(112, 10).Value = 1512.4667
$ws.Cells.Item(112, 12).Value = 4537.4001
$ws.Cells.Item(112, 14).Value = -6753.4001
$ws.Cells.Item(135, 8).Value = 2332.2222
$ws.Cells.Item(135, 9).Value = 1232
$ws.Cells.Item(135, 11).Value = 11088
$ws.Cells.Item(135, 13).Value = -8553
$ws.Cells.Item(137, 8).Value = 1196.8572
$ws.Cells.Item(137, 9).Value = 951.37933
$ws.Cells.Item(137, 11).Value = 2854.13799
$ws.Cells.Item(137, 13).Value = -304.1379900000002
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1474.2222
$ws.Cells.Item(45, 9).Value = 1012.8
$ws.Cells.Item(45, 10).Value = 2051
$ws.Cells.Item(45, 11).Value = 1012.8
$ws.Cells.Item(45, 12).Value = 2051
$ws.Cells.Item(45, 13).Value = -635.8
$ws.Cells.Item(45, 14).Value = -2805
$ws.Cells.Item(61, 8).Value = 1059.7778
$ws.Cells.Item(61, 9).Value = 704.8
$ws.Cells.Item(61, 10).Value = 1503.5
$ws.Cells.Item(61, 11).Value = 704.8
$ws.Cells.Item(61, 12).Value = 1503.5
$ws.Cells.Item(61, 13).Value = -492.8
$ws.Cells.Item(61, 14).Value = -1927.5
$ws.Cells.Item(74, 8).Value = 1228.9412
$ws.Cells.Item(74, 9).Value = 1255.75
$ws.Cells.Item(74, 11).Value = 1255.75
$ws.Cells.Item(74, 13).Value = -381.75
$ws.Cells.Item(77, 8).Value = 1228.9412
$ws.Cells.Item(77, 9).Value = 1255.75
$ws.Cells.Item(77, 11).Value = 6278.75
$ws.Cells.Item(77, 13).Value = -1910.75
$ws.Cells.Item(102, 8).Value = 4245.5713
$ws.Cells.Item(102, 9).Value = 4640
$ws.Cells.Item(102, 11).Value = 4640
$ws.Cells.Item(102, 13).Value = -3018
$ws.Cells.Item(132, 8).Value = 2247.8572
$ws.Cells.Item(132, 9).Value = 1824.0416
$ws.Cells.Item(132, 10).Value = 3172.5454
$ws.Cells.Item(132, 11).Value = 5472.1248
$ws.Cells.Item(132, 12).Value = 9517.636200000001
$ws.Cells.Item(132, 13).Value = -2942.1248
$ws.Cells.Item(132, 14).Value = -14577.6362
$ws.Cells.Item(136, 8).Value = 1059.7778
$ws.Cells.Item(136, 9).Value = 704.8
$ws.Cells.Item(136, 10).Value = 1503.5
$ws.Cells.Item(136, 11).Value = 2114.4
$ws.Cells.Item(136, 12).Value = 4510.5
$ws.Cells.Item(136, 13).Value = 435.6000000000004
$ws.Cells.Item(136, 14).Value = -9610.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 6866.478
$ws.Cells.Item(107, 9).Value = 1142.2273
$ws.Cells.Item(107, 10).Value = 132800
$ws.Cells.Item(107, 11).Value = 1142.2273
$ws.Cells.Item(107, 12).Value = 132800
$ws.Cells.Item(107, 13).Value = 777.7727
$ws.Cells.Item(107, 14).Value = -136640
$ws.Cells.Item(134, 8).Value = 31229.686
$ws.Cells.Item(134, 9).Value = 2541.652
$ws.Cells.Item(134, 10).Value = 86215.086
$ws.Cells.Item(134, 11).Value = 7624.956
$ws.Cells.Item(134, 12).Value = 258645.258
$ws.Cells.Item(134, 13).Value = -5089.956
$ws.Cells.Item(134, 14).Value = -263715.258
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 2000
$ws.Cells.Item(4, 10).Value = 2000
$ws.Cells.Item(4, 12).Value = 2000
$ws.Cells.Item(4, 13).Value = -2224
$ws.Cells.Item(31, 8).Value = 2352.0605
$ws.Cells.Item(31, 9).Value = 2253.2183
$ws.Cells.Item(31, 10).Value = 2846.2727
$ws.Cells.Item(31, 11).Value = 2253.2183
$ws.Cells.Item(31, 12).Value = 2846.2727
$ws.Cells.Item(31, 13).Value = -1958.2183
$ws.Cells.Item(31, 14).Value = -3436.2727
$ws.Cells.Item(34, 8).Value = 2352.0605
$ws.Cells.Item(34, 9).Value = 2253.2183
$ws.Cells.Item(34, 10).Value = 2846.2727
$ws.Cells.Item(34, 11).Value = 2253.2183
$ws.Cells.Item(34, 12).Value = 2846.2727
$ws.Cells.Item(34, 13).Value = -2051.2183
$ws.Cells.Item(34, 14).Value = -3250.2727
$ws.Cells.Item(58, 8).Value = 3084.8696
$ws.Cells.Item(58, 9).Value = 657.2069
$ws.Cells.Item(58, 10).Value = 7226.1763
$ws.Cells.Item(58, 11).Value = 657.2069
$ws.Cells.Item(58, 12).Value = 7226.1763
$ws.Cells.Item(58, 13).Value = -454.2069
$ws.Cells.Item(58, 14).Value = -7632.1763
$ws.Cells.Item(132, 8).Value = 2212.4092
$ws.Cells.Item(132, 9).Value = 1939.6471
$ws.Cells.Item(132, 10).Value = 3139.8
$ws.Cells.Item(132, 11).Value = 5818.9413
$ws.Cells.Item(132, 12).Value = 9419.400000000001
$ws.Cells.Item(132, 13).Value = -3288.9413
$ws.Cells.Item(132, 14).Value = -14479.4
$ws.Cells.Item(134, 8).Value = 1729.579
$ws.Cells.Item(134, 9).Value = 769.2727
$ws.Cells.Item(134, 10).Value = 3050
$ws.Cells.Item(134, 11).Value = 2307.8181
$ws.Cells.Item(134, 12).Value = 9150
$ws.Cells.Item(134, 13).Value = 227.1819
$ws.Cells.Item(134, 14).Value = -14220
$ws.Cells.Item(136, 8).Value = 3084.8696
$ws.Cells.Item(136, 9).Value = 657.2069
$ws.Cells.Item(136, 10).Value = 7226.1763
$ws.Cells.Item(136, 11).Value = 1971.6207
$ws.Cells.Item(136, 12).Value = 21678.5289
$ws.Cells.Item(136, 13).Value = 578.3793000000001
$ws.Cells.Item(136, 14).Value = -26778.5289
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(109, 8).Value = 45339.26
$ws.Cells.Item(109, 9).Value = 91563.91
$ws.Cells.Item(109, 10).Value = 2966.6667
$ws.Cells.Item(109, 11).Value = 274691.73
$ws.Cells.Item(109, 12).Value = 8900.000100000001
$ws.Cells.Item(109, 13).Value = -273651.73
$ws.Cells.Item(109, 14).Value = -10980.0001
$ws.Cells.Item(113, 8).Value = 726.36365
$ws.Cells.Item(113, 9).Value = 1000
$ws.Cells.Item(113, 10).Value = 623.75
$ws.Cells.Item(113, 11).Value = 3000
$ws.Cells.Item(113, 12).Value = 1871.25
$ws.Cells.Item(113, 13).Value = -830
$ws.Cells.Item(113, 14).Value = -6211.25
$ws.Cells.Item(122, 8).Value = 1112700.4
$ws.Cells.Item(122, 9).Value = 652
$ws.Cells.Item(122, 10).Value = 1430428.4
$ws.Cells.Item(122, 11).Value = 5868
$ws.Cells.Item(122, 12).Value = 12873855.6
$ws.Cells.Item(122, 13).Value = -3418
$ws.Cells.Item(122, 14).Value = -12878755.6
$ws.Cells.Item(123, 8).Value = 2922.8572
$ws.Cells.Item(123, 9).Value = 780
$ws.Cells.Item(123, 10).Value = 5780
$ws.Cells.Item(123, 11).Value = 2340
$ws.Cells.Item(123, 12).Value = 17340
$ws.Cells.Item(123, 13).Value = 110
$ws.Cells.Item(123, 14).Value = -22240
$ws.Cells.Item(125, 8).Value = 2589.8572
$ws.Cells.Item(125, 9).Value = 1810
$ws.Cells.Item(125, 10).Value = 3174.75
$ws.Cells.Item(125, 11).Value = 5430
$ws.Cells.Item(125, 12).Value = 9524.25
$ws.Cells.Item(125, 13).Value = -510
$ws.Cells.Item(125, 14).Value = -19364.25
$ws.Cells.Item(129, 8).Value = 2102.875
$ws.Cells.Item(129, 10).Value = 3311
$ws.Cells.Item(129, 12).Value = 9933
$ws.Cells.Item(129, 14).Value = -19933
$ws.Cells.Item(131, 8).Value = 7157335.5
$ws.Cells.Item(131, 9).Value = 45917824
$ws.Cells.Item(131, 10).Value = 1552.8616
$ws.Cells.Item(131, 11).Value = 137753472
$ws.Cells.Item(131, 12).Value = 4658.5848
$ws.Cells.Item(131, 13).Value = -137748432
$ws.Cells.Item(131, 14).Value = -14738.5848
$ws.Cells.Item(132, 8).Value = 1360
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 10).Value = 1360
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 12).ClearContents()
$ws.Cells.Item(132, 13).Value = 12240
$ws.Cells.Item(132, 14).Value = -17300
$ws.Cells.Item(137, 8).Value = 56113.15
$ws.Cells.Item(137, 9).Value = 2373.5715
$ws.Cells.Item(137, 10).Value = 181505.5
$ws.Cells.Item(137, 11).Value = 7120.7145
$ws.Cells.Item(137, 12).Value = 544516.5
$ws.Cells.Item(137, 13).Value = -2020.7145
$ws.Cells.Item(137, 14).Value = -554716.5
$ws.Cells.Item(138, 8).Value = 1295.4546
$ws.Cells.Item(138, 9).Value = 925
$ws.Cells.Item(138, 11).Value = 2775
$ws.Cells.Item(138, 13).Value = 2365
$ws.Cells.Item(139, 8).Value = 38336.297
$ws.Cells.Item(139, 9).Value = 42503.332
$ws.Cells.Item(139, 11).Value = 127509.996
$ws.Cells.Item(139, 13).Value = -122369.996
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4213.829
$ws.Cells.Item(70, 9).Value = 4149.8887
$ws.Cells.Item(70, 10).Value = 4337.143
$ws.Cells.Item(70, 11).Value = 4149.8887
$ws.Cells.Item(70, 12).Value = 4337.143
$ws.Cells.Item(70, 13).Value = -3879.8887
$ws.Cells.Item(70, 14).Value = -4877.143
$ws.Cells.Item(73, 8).Value = 4213.829
$ws.Cells.Item(73, 9).Value = 4149.8887
$ws.Cells.Item(73, 10).Value = 4337.143
$ws.Cells.Item(73, 11).Value = 4149.8887
$ws.Cells.Item(73, 12).Value = 4337.143
$ws.Cells.Item(73, 13).Value = -3213.8887
$ws.Cells.Item(73, 14).Value = -6209.143
$ws.Cells.Item(132, 8).Value = 3754.6365
$ws.Cells.Item(132, 9).Value = 3489.111
$ws.Cells.Item(132, 11).Value = 10467.333
$ws.Cells.Item(132, 13).Value = -7937.332999999999
$ws.Cells.Item(133, 8).Value = 44181.816
$ws.Cells.Item(133, 10).Value = 44181.816
$ws.Cells.Item(133, 12).Value = 44181.816
$ws.Cells.Item(133, 14).Value = -54301.816
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 22558.8
$ws.Cells.Item(61, 9).Value = 22558.8
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 11).Value = 22558.8
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 13).ClearContents()
$ws.Cells.Item(61, 14).Value = -22356.8
$ws.Cells.Item(113, 8).Value = 22558.8
$ws.Cells.Item(113, 9).Value = 22558.8
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 22558.8
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).ClearContents()
$ws.Cells.Item(113, 14).Value = -20388.8
$ws.Cells.Item(136, 8).Value = 3473.9
$ws.Cells.Item(136, 9).Value = 1673.1852
$ws.Cells.Item(136, 10).Value = 7213.846
$ws.Cells.Item(136, 11).Value = 5019.5556
$ws.Cells.Item(136, 12).Value = 21641.538
$ws.Cells.Item(136, 13).Value = -2469.5556
$ws.Cells.Item(136, 14).Value = -26741.538
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(5, 8).Value = 2286815.8
$ws.Cells.Item(5, 9).Value = 805.5
$ws.Cells.Item(5, 10).Value = 3201220
$ws.Cells.Item(5, 11).Value = 805.5
$ws.Cells.Item(5, 12).Value = 3201220
$ws.Cells.Item(5, 13).Value = -693.5
$ws.Cells.Item(5, 14).Value = -3201444
$ws.Cells.Item(122, 8).Value = 464.66666
$ws.Cells.Item(122, 9).Value = 464.66666
$ws.Cells.Item(122, 11).Value = 1393.99998
$ws.Cells.Item(122, 13).Value = 1056.00002
$ws.Cells.Item(132, 8).Value = 1547.6863
$ws.Cells.Item(132, 9).Value = 1535.9445
$ws.Cells.Item(132, 10).Value = 1554.091
$ws.Cells.Item(132, 11).Value = 4607.833500000001
$ws.Cells.Item(132, 12).Value = 4662.272999999999
$ws.Cells.Item(132, 13).Value = -2077.833500000001
$ws.Cells.Item(132, 14).Value = -9722.272999999999
$ws.Cells.Item(140, 8).Value = 49292.11
$ws.Cells.Item(140, 10).Value = 49292.11
$ws.Cells.Item(140, 12).Value = 49292.11
$ws.Cells.Item(140, 14).Value = -59652.11
